$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated probabilities (column C)
$ws.Range("C3").Value = 0.44
$ws.Range("C5").Value = 0.08
$ws.Range("C12").Value = 0.03
$ws.Range("C13").Value = 0.03

# Updated generation ranges (column E).
# Order matches the shared-string insertion order recorded in the target file.
$ws.Range("E3").Value = "1-44"
$ws.Range("E13").Value = "98-100"
$ws.Range("E4").Value = "45-54"
$ws.Range("E5").Value = "55-62"
$ws.Range("E6").Value = "63-68"
$ws.Range("E7").Value = "69-78"
$ws.Range("E8").Value = "79-83"
$ws.Range("E9").Value = "84-87"
$ws.Range("E10").Value = "88-91"
$ws.Range("E11").Value = "92-94"
$ws.Range("E12").Value = "95-97"

# Update the active selection to G9
$ws.Range("G9").Select()
